$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OOMII_DB")

# Before this edit, row 16 (the last data row) spans the sheet's full width
# (columns A:OR) with an explicit - but empty - cell for every unused
# column. The edit:
#   1. Strips those empty placeholder cells from row 16, leaving only the
#      columns that actually hold data (matching how every other data row
#      in the sheet, 2-15, is authored: sparse, values-only).
#   2. Appends a new row 17 with the same data as row 16 used to have,
#      written the same sparse way.
#   3. Appends a new row 18 that is a full copy of row 16 (including the
#      empty placeholder cells for every unused column across the sheet's
#      full width) - i.e. row 18 ends up looking like row 16 did originally.

$maxCol = 408

# Capture row 16's current text per column before we start mutating it.
$row16Text = @{}
for ($c = 1; $c -le $maxCol; $c++) {
    $row16Text[$c] = $ws.Cells.Item(16, $c).Text
}

# Step 3 (do this first, while row 16 still has its original full-width
# layout): copy row 16 verbatim into row 18. A range copy - unlike direct
# value assignment - preserves an explicit empty cell node for every
# column, reproducing row 16's original look.
$srcRow = $ws.Range("A16:OR16")
$dstRow = $ws.Range("A18:OR18")
$srcRow.Copy($dstRow)

# Step 2: write row 17 cell-by-cell, sparse (only columns that actually
# held a value get a cell).
for ($c = 1; $c -le $maxCol; $c++) {
    $srcText = $row16Text[$c]
    if ($srcText.Length -gt 0) {
        $dst = $ws.Cells.Item(17, $c)
        # Force text storage so numeric-looking values (e.g. "0.0625",
        # "94", "-0") round-trip as strings instead of being reinterpreted
        # as numbers.
        $dst.NumberFormat = "@"
        $dst.Value = $srcText
    }
}

# Step 1: strip the empty placeholder cells from row 16 itself, leaving it
# sparse like the rest of the data rows.
for ($c = 1; $c -le $maxCol; $c++) {
    $srcText = $row16Text[$c]
    if ($srcText.Length -eq 0) {
        $ws.Cells.Item(16, $c).ClearContents()
    }
}
